$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.288.57'
$ws.Range('E2').Value = '  -0.34%  '
$ws.Range('D3').Value = '2.271.11'
$ws.Range('E3').Value = '  -0.44%  '
$ws.Range('E4').Value = '  -0.08%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '308.63'
$ws.Range('E5').Value = '  +0.44%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '97.53'
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('E7').Value = '  -1.03%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('E9').Value = '  -1.40%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.98'
$ws.Range('E10').Value = '  -3.00%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0809'
$ws.Range('E11').Value = '  +0.76%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '6.84'
$ws.Range('E13').Value = '  +1.79%  '
$ws.Range('D14').Value = '2.623.51'
$ws.Range('E14').Value = '  -0.44%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '14.58'
$ws.Range('E15').Value = '  +0.56%  '
$ws.Range('D16').Value = '2.261.79'
$ws.Range('E16').Value = '  -1.07%  '
$ws.Range('E17').Value = '  -1.74%  '
$ws.Range('D18').Value = '42.169.12'
$ws.Range('E18').Value = '  -0.36%  '
$ws.Range('E19').Value = '  -2.74%  '
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('E20').Value = '  -1.07%  '
$ws.Range('E21').Value = '  -1.10%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '67.59'
$ws.Range('E22').Value = '  -0.33%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '236.59'
$ws.Range('E23').Value = '  -2.26%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.60'
$ws.Range('E24').Value = '  -0.11%  '
$ws.Range('E25').Value = '  +1.15%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '23.59'
$ws.Range('E27').Value = '  -1.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '37.14'
$ws.Range('E28').Value = '  -2.25%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  +0.50%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '163.20'
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('E32').Value = '  -0.35%  '
$ws.Range('E33').Value = '  +0.00%  '
$ws.Range('E34').Value = '  -1.07%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '17.58'
$ws.Range('E35').Value = '  +2.13%  '
$ws.Range('E36').Value = '  -2.09%  '
$ws.Range('E37').Value = '  -0.30%  '
$ws.Range('E38').Value = '  -2.78%  '
$ws.Range('E39').Value = '  -0.31%  '
$ws.Range('E40').Value = '  -2.02%  '
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('E42').Value = '  -5.22%  '
$ws.Range('D43').Value = '1.946.22'
$ws.Range('E43').Value = '  -2.85%  '
$ws.Range('E44').Value = '  -1.47%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '18.88'
$ws.Range('E45').Value = '  -0.87%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '2.94'
$ws.Range('E46').Value = '  -2.28%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '9.77'
$ws.Range('E47').Value = '  -3.03%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '54.36'
$ws.Range('E48').Value = '  +1.60%  '
$ws.Range('D49').Value = '2.495.53'
$ws.Range('E49').Value = '  -0.44%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '91.97'
$ws.Range('E50').Value = '  -1.21%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '71.56'
$ws.Range('E51').Value = '  -1.18%  '
